$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.356.62'
$ws.Range("E2").Value = '  -4.86%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.260.13'
$ws.Range("E3").Value = '  -7.67%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.33'
$ws.Range("E5").Value = '  -5.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.89'
$ws.Range("E6").Value = '  -12.38%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.251.13'
$ws.Range("E8").Value = '  -7.87%  '
$ws.Range("E9").Value = '  -10.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.173'
$ws.Range("E10").Value = '  -13.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.85'
$ws.Range("E11").Value = '  -3.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.509'
$ws.Range("E12").Value = '  -13.24%  '
$ws.Range("E13").Value = '  -17.00%  '
$ws.Range("E14").Value = '  -11.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.776.31'
$ws.Range("E15").Value = '  -7.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.368.20'
$ws.Range("E16").Value = '  -4.93%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '549.99'
$ws.Range("E17").Value = '  -9.53%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.257.41'
$ws.Range("E18").Value = '  -7.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.30'
$ws.Range("E19").Value = '  -13.36%  '
$ws.Range("E20").Value = '  -5.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.24'
$ws.Range("E21").Value = '  -14.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.771'
$ws.Range("E22").Value = '  -12.93%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.95'
$ws.Range("E23").Value = '  -12.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.88'
$ws.Range("E24").Value = '  -12.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.66'
$ws.Range("E25").Value = '  -12.84%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").Value = '  -14.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.18'
$ws.Range("E28").Value = '  -10.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '29.60'
$ws.Range("E29").Value = '  -12.50%  '
$ws.Range("E30").Value = '  -16.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.68'
$ws.Range("E31").Value = '  -11.95%  '
$ws.Range("E32").Value = '  -11.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '548.53'
$ws.Range("E33").Value = '  -12.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.69'
$ws.Range("E34").Value = '  -17.77%  '
$ws.Range("E35").Value = '  -14.89%  '
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0449'
$ws.Range("E37").Value = '  -5.65%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '53.66'
$ws.Range("E38").Value = '  -5.59%  '
$ws.Range("E39").Value = '  -13.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.27'
$ws.Range("E40").Value = '  -14.45%  '
$ws.Range("E41").Value = '  -11.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.947.78'
$ws.Range("E42").Value = '  -12.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.67'
$ws.Range("E43").Value = '  -23.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.264'
$ws.Range("E44").Value = '  -15.56%  '
$ws.Range("E45").Value = '  -19.98%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '26.61'
$ws.Range("E46").Value = '  -17.07%  '
$ws.Range("E47").Value = '  -19.87%  '
$ws.Range("E48").Value = '  -15.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '127.15'
$ws.Range("E50").Value = '  -4.44%  '
$ws.Range("E51").Value = '  -12.10%  '
